$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these Price cells as plain text (they are stored as text in the workbook,
# e.g. "148.80"/"6.60"/"1.60" - without this, Excel would coerce them to numbers
# and silently drop the significant trailing zero).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "59.912.50"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "2.660.07"
$ws.Range("E3").Value = "  +2.33%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "520.93"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D6").Value = "148.80"
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "0.574"
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "2.690.75"
$ws.Range("E9").Value = "  +3.36%  "
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("D11").Value = "0.107"
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("D14").Value = "3.132.71"
$ws.Range("E14").Value = "  +2.32%  "
$ws.Range("D15").Value = "59.791.50"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "21.51"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("E17").Value = "  +1.93%  "
$ws.Range("D18").Value = "2.694.07"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("D19").Value = "4.65"
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "348.66"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").Value = "10.65"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("E22").Value = "  +3.49%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "61.25"
$ws.Range("E24").Value = "  +1.72%  "
$ws.Range("D25").Value = "0.428"
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("D26").Value = "2.780.11"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").Value = "0.162"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "0.991"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").Value = "0.0₃0833"
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("E30").Value = "  +3.29%  "
$ws.Range("D31").Value = "6.60"
$ws.Range("E31").Value = "  +10.90%  "
$ws.Range("E32").Value = "  -0.35%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "19.16"
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("B34").Value = "PancakeSwap"
$ws.Range("C34").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D34").Value = "1.60"
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("D35").Value = "1.09"
$ws.Range("E35").Value = "  +20.93%  "
$ws.Range("D36").Value = "148.48"
$ws.Range("E36").Value = "  -0.91%  "
$ws.Range("E37").Value = "  +3.08%  "
$ws.Range("E38").Value = "  +3.49%  "
$ws.Range("D39").Value = "0.883"
$ws.Range("E39").Value = "  +2.83%  "
$ws.Range("D40").Value = "36.66"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("D41").Value = "3.75"
$ws.Range("E41").Value = "  +3.80%  "
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("D43").Value = "288.19"
$ws.Range("E43").Value = "  +1.16%  "
$ws.Range("D44").Value = "0.628"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("D45").Value = "0.0999"
$ws.Range("E45").Value = "  -0.10%  "
$ws.Range("D46").Value = "0.993"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").Value = "19.88"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").Value = "4.81"
$ws.Range("E49").Value = "  +2.72%  "
$ws.Range("E50").Value = "  +0.89%  "
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.999.61"
$ws.Range("E51").Value = "  +3.00%  "
